$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 615-616, pushing the existing rows 615..702 down to 617..704.
$ws.Rows("615:616").Insert()

# New row 615: Asterix, "1a (guarda)" (quality unchanged from the row that used to sit here),
# date 2023-08-04 (serial 45142), volume/price/origin updated.
$ws.Cells.Item(615, 1).Value2 = 7
$ws.Cells.Item(615, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(615, 3).Value2 = "Ñuble"
$ws.Cells.Item(615, 4).Value2 = 45142
$ws.Cells.Item(615, 5).Value2 = 16
$ws.Cells.Item(615, 6).Value2 = 100114001
$ws.Cells.Item(615, 7).Value2 = "Papa"
$ws.Cells.Item(615, 8).Value2 = "Asterix"
$ws.Cells.Item(615, 9).Value2 = "1a (guarda)"
$ws.Cells.Item(615, 10).Value2 = 150
$ws.Cells.Item(615, 11).Value2 = 18000
$ws.Cells.Item(615, 12).Value2 = 18000
$ws.Cells.Item(615, 13).Value2 = 18000
$ws.Cells.Item(615, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(615, 15).Value2 = "Región de Los Lagos"
$ws.Cells.Item(615, 16).Value2 = 720
$ws.Cells.Item(615, 17).Value2 = 25
$ws.Cells.Item(615, 18).Value2 = "Hortaliza"

# New row 616: Asterix, "2a (guarda)", same date 2023-08-04 (serial 45142).
$ws.Cells.Item(616, 1).Value2 = 7
$ws.Cells.Item(616, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(616, 3).Value2 = "Ñuble"
$ws.Cells.Item(616, 4).Value2 = 45142
$ws.Cells.Item(616, 5).Value2 = 16
$ws.Cells.Item(616, 6).Value2 = 100114001
$ws.Cells.Item(616, 7).Value2 = "Papa"
$ws.Cells.Item(616, 8).Value2 = "Asterix"
$ws.Cells.Item(616, 9).Value2 = "2a (guarda)"
$ws.Cells.Item(616, 10).Value2 = 150
$ws.Cells.Item(616, 11).Value2 = 16000
$ws.Cells.Item(616, 12).Value2 = 16000
$ws.Cells.Item(616, 13).Value2 = 16000
$ws.Cells.Item(616, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(616, 15).Value2 = "Región de Los Lagos"
$ws.Cells.Item(616, 16).Value2 = 640
$ws.Cells.Item(616, 17).Value2 = 25
$ws.Cells.Item(616, 18).Value2 = "Hortaliza"
